# Generate Report for Handoff
# Replace the old GUID-based file name ("01c933fd-2d8c-45b8-9702-90196683f065")
# with the new one ("afb0c21f-1429-4d00-ae16-26a10c9a1f23") throughout the
# workbook, and bump the associated handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "01c933fd-2d8c-45b8-9702-90196683f065"
$newGuid = "afb0c21f-1429-4d00-ae16-26a10c9a1f23"

$oldHash = "c50ee5e3781b66ee5bd5d0e0eba028f70308b87a"
$newHash = "46890682b58ddfdbeb45b52fe03d5d5b7ff4d997"

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-23 15:12:50"

# --- Sheet "zh-cn" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-23 15:12:35"

# --- Sheet "de-de" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-23 15:12:50"
